$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two now-obsolete data rows (previously rows 4 and 5);
# this shifts nothing below them, it just drops the old branch entries.
$ws.Rows("4:5").Delete() | Out-Null

# Row 2: new CNPJ / Fantasia data (Try/Except style placeholder test values)
$ws.Range("A2").Value = "88.543.679/0001-06"

# Row 3: second new CNPJ / Fantasia entry
$ws.Range("A3").Value = "46.346.081/0001-87"
$ws.Range("B3").Value = "TesteCNPJ01"
$ws.Range("C3").Value = "TesteCNPJ01"
$ws.Range("B2").Value = "TesteCNPJ00"
$ws.Range("C2").Value = "TesteCNPJ00"

# Mark the CNPJ cell of the first data row with an underline style
$ws.Range("A2").Font.Underline = $true

# Leave a styled-but-empty marker cell further down the sheet
$ws.Range("C8").Font.Underline = $true

# Widen the Fantasia/Nome columns to fit the new text
$ws.Columns("B:C").ColumnWidth = 45.6

# Match the selection left behind in the saved file
$ws.Range("C8").Select() | Out-Null
